$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 <- original row 36
$ws.Cells.Item(2, 4).Value = 44540
$ws.Cells.Item(2, 9).Value = "Primera"
$ws.Cells.Item(2, 10).Value = 300
$ws.Cells.Item(2, 11).Value = 900
$ws.Cells.Item(2, 12).Value = 1000
$ws.Cells.Item(2, 13).Value = 950
$ws.Cells.Item(2, 14).Value = "`$/atado 1,5 a 2 kilos"
$ws.Cells.Item(2, 16).Value = 475
$ws.Cells.Item(2, 17).Value = 2

# Row 3 <- original row 14
$ws.Cells.Item(3, 4).Value = 44363
$ws.Cells.Item(3, 9).Value = "Primera"
$ws.Cells.Item(3, 10).Value = 250
$ws.Cells.Item(3, 11).Value = 2500
$ws.Cells.Item(3, 12).Value = 2800
$ws.Cells.Item(3, 13).Value = 2650
$ws.Cells.Item(3, 14).Value = "`$/atado 1,5 a 2 kilos"
$ws.Cells.Item(3, 16).Value = 1325
$ws.Cells.Item(3, 17).Value = 2

# Row 4 <- original row 37
$ws.Cells.Item(4, 4).Value = 44817
$ws.Cells.Item(4, 9).Value = "Primera"
$ws.Cells.Item(4, 10).Value = 300
$ws.Cells.Item(4, 11).Value = 900
$ws.Cells.Item(4, 12).Value = 1000
$ws.Cells.Item(4, 13).Value = 950
$ws.Cells.Item(4, 14).Value = "`$/atado 1,5 a 2 kilos"
$ws.Cells.Item(4, 16).Value = 475
$ws.Cells.Item(4, 17).Value = 2

# Row 5 <- original row 30
$ws.Cells.Item(5, 4).Value = 44291
$ws.Cells.Item(5, 9).Value = "Primera"
$ws.Cells.Item(5, 10).Value = 250
$ws.Cells.Item(5, 11).Value = 1800
$ws.Cells.Item(5, 12).Value = 2000
$ws.Cells.Item(5, 13).Value = 1900
$ws.Cells.Item(5, 14).Value = "`$/atado 1,5 a 2 kilos"
$ws.Cells.Item(5, 16).Value = 950
$ws.Cells.Item(5, 17).Value = 2

# Row 6 <- original row 19
$ws.Cells.Item(6, 4).Value = 44572
$ws.Cells.Item(6, 9).Value = "Primera"
$ws.Cells.Item(6, 10).Value = 300
$ws.Cells.Item(6, 11).Value = 1400
$ws.Cells.Item(6, 12).Value = 1500
$ws.Cells.Item(6, 13).Value = 1450
$ws.Cells.Item(6, 14).Value = "`$/atado 1,5 a 2 kilos"
$ws.Cells.Item(6, 16).Value = 725
$ws.Cells.Item(6, 17).Value = 2

# Row 7 <- original row 34
$ws.Cells.Item(7, 4).Value = 45070
$ws.Cells.Item(7, 9).Value = "Primera"
$ws.Cells.Item(7, 10).Value = 270
$ws.Cells.Item(7, 11).Value = 1000
$ws.Cells.Item(7, 12).Value = 1500
$ws.Cells.Item(7, 13).Value = 1250
$ws.Cells.Item(7, 14).Value = "`$/atado 1,5 a 2 kilos"
$ws.Cells.Item(7, 16).Value = 625
$ws.Cells.Item(7, 17).Value = 2

# Row 8 <- original row 29
$ws.Cells.Item(8, 4).Value = 44789
$ws.Cells.Item(8, 9).Value = "Primera"
$ws.Cells.Item(8, 10).Value = 300
$ws.Cells.Item(8, 11).Value = 1400
$ws.Cells.Item(8, 12).Value = 1500
$ws.Cells.Item(8, 13).Value = 1450
$ws.Cells.Item(8, 14).Value = "`$/atado 1,5 a 2 kilos"
$ws.Cells.Item(8, 16).Value = 725
$ws.Cells.Item(8, 17).Value = 2

# Row 9 <- original row 16
$ws.Cells.Item(9, 4).Value = 44302
$ws.Cells.Item(9, 9).Value = "Primera"
$ws.Cells.Item(9, 10).Value = 300
$ws.Cells.Item(9, 11).Value = 900
$ws.Cells.Item(9, 12).Value = 1000
$ws.Cells.Item(9, 13).Value = 950
$ws.Cells.Item(9, 14).Value = "`$/atado 1,5 a 2 kilos"
$ws.Cells.Item(9, 16).Value = 475
$ws.Cells.Item(9, 17).Value = 2

# Row 10 <- original row 33
$ws.Cells.Item(10, 4).Value = 44972
$ws.Cells.Item(10, 9).Value = "Primera"
$ws.Cells.Item(10, 10).Value = 350
$ws.Cells.Item(10, 11).Value = 800
$ws.Cells.Item(10, 12).Value = 1000
$ws.Cells.Item(10, 13).Value = 943
$ws.Cells.Item(10, 14).Value = "`$/atado 1,5 a 2 kilos"
$ws.Cells.Item(10, 16).Value = 472
$ws.Cells.Item(10, 17).Value = 2

# Row 11 <- original row 25
$ws.Cells.Item(11, 4).Value = 44202
$ws.Cells.Item(11, 9).Value = "Primera"
$ws.Cells.Item(11, 10).Value = 250
$ws.Cells.Item(11, 11).Value = 1800
$ws.Cells.Item(11, 12).Value = 2000
$ws.Cells.Item(11, 13).Value = 1900
$ws.Cells.Item(11, 14).Value = "`$/atado 1,5 a 2 kilos"
$ws.Cells.Item(11, 16).Value = 950
$ws.Cells.Item(11, 17).Value = 2

# Row 12 <- original row 35
$ws.Cells.Item(12, 4).Value = 45062
$ws.Cells.Item(12, 9).Value = "Primera"
$ws.Cells.Item(12, 10).Value = 380
$ws.Cells.Item(12, 11).Value = 1800
$ws.Cells.Item(12, 12).Value = 2000
$ws.Cells.Item(12, 13).Value = 1895
$ws.Cells.Item(12, 14).Value = "`$/atado 1,5 a 2 kilos"
$ws.Cells.Item(12, 16).Value = 948
$ws.Cells.Item(12, 17).Value = 2

# Row 13 <- original row 44
$ws.Cells.Item(13, 4).Value = 44253
$ws.Cells.Item(13, 9).Value = "Primera"
$ws.Cells.Item(13, 10).Value = 250
$ws.Cells.Item(13, 11).Value = 1800
$ws.Cells.Item(13, 12).Value = 2000
$ws.Cells.Item(13, 13).Value = 1900
$ws.Cells.Item(13, 14).Value = "`$/atado 1,5 a 2 kilos"
$ws.Cells.Item(13, 16).Value = 950
$ws.Cells.Item(13, 17).Value = 2

# Row 14 <- original row 39
$ws.Cells.Item(14, 4).Value = 44525
$ws.Cells.Item(14, 9).Value = "Primera"
$ws.Cells.Item(14, 10).Value = 300
$ws.Cells.Item(14, 11).Value = 1400
$ws.Cells.Item(14, 12).Value = 1500
$ws.Cells.Item(14, 13).Value = 1450
$ws.Cells.Item(14, 14).Value = "`$/atado 1,5 a 2 kilos"
$ws.Cells.Item(14, 16).Value = 725
$ws.Cells.Item(14, 17).Value = 2

# Row 15 <- original row 3
$ws.Cells.Item(15, 4).Value = 44435
$ws.Cells.Item(15, 9).Value = "Primera"
$ws.Cells.Item(15, 10).Value = 300
$ws.Cells.Item(15, 11).Value = 900
$ws.Cells.Item(15, 12).Value = 1000
$ws.Cells.Item(15, 13).Value = 950
$ws.Cells.Item(15, 14).Value = "`$/atado 1,5 a 2 kilos"
$ws.Cells.Item(15, 16).Value = 475
$ws.Cells.Item(15, 17).Value = 2

# Row 16 <- original row 27
$ws.Cells.Item(16, 4).Value = 44229
$ws.Cells.Item(16, 9).Value = "Primera"
$ws.Cells.Item(16, 10).Value = 250
$ws.Cells.Item(16, 11).Value = 1800
$ws.Cells.Item(16, 12).Value = 2000
$ws.Cells.Item(16, 13).Value = 1900
$ws.Cells.Item(16, 14).Value = "`$/atado 1,5 a 2 kilos"
$ws.Cells.Item(16, 16).Value = 950
$ws.Cells.Item(16, 17).Value = 2

# Row 17 <- original row 10
$ws.Cells.Item(17, 4).Value = 44726
$ws.Cells.Item(17, 9).Value = "Primera"
$ws.Cells.Item(17, 10).Value = 250
$ws.Cells.Item(17, 11).Value = 2500
$ws.Cells.Item(17, 12).Value = 2800
$ws.Cells.Item(17, 13).Value = 2650
$ws.Cells.Item(17, 14).Value = "`$/atado 1,5 a 2 kilos"
$ws.Cells.Item(17, 16).Value = 1325
$ws.Cells.Item(17, 17).Value = 2

# Row 18 <- original row 13
$ws.Cells.Item(18, 4).Value = 44266
$ws.Cells.Item(18, 9).Value = "Primera"
$ws.Cells.Item(18, 10).Value = 300
$ws.Cells.Item(18, 11).Value = 1700
$ws.Cells.Item(18, 12).Value = 1800
$ws.Cells.Item(18, 13).Value = 1750
$ws.Cells.Item(18, 14).Value = "`$/atado 1,5 a 2 kilos"
$ws.Cells.Item(18, 16).Value = 875
$ws.Cells.Item(18, 17).Value = 2

# Row 19 <- original row 23
$ws.Cells.Item(19, 4).Value = 44795
$ws.Cells.Item(19, 9).Value = "Primera"
$ws.Cells.Item(19, 10).Value = 250
$ws.Cells.Item(19, 11).Value = 1800
$ws.Cells.Item(19, 12).Value = 2000
$ws.Cells.Item(19, 13).Value = 1900
$ws.Cells.Item(19, 14).Value = "`$/atado 1,5 a 2 kilos"
$ws.Cells.Item(19, 16).Value = 950
$ws.Cells.Item(19, 17).Value = 2

# Row 20 <- original row 9
$ws.Cells.Item(20, 4).Value = 44161
$ws.Cells.Item(20, 9).Value = "Primera"
$ws.Cells.Item(20, 10).Value = 270
$ws.Cells.Item(20, 11).Value = 900
$ws.Cells.Item(20, 12).Value = 1000
$ws.Cells.Item(20, 13).Value = 950
$ws.Cells.Item(20, 14).Value = "`$/atado 1,5 a 2 kilos"
$ws.Cells.Item(20, 16).Value = 475
$ws.Cells.Item(20, 17).Value = 2

# Row 21 <- original row 18
$ws.Cells.Item(21, 4).Value = 44390
$ws.Cells.Item(21, 9).Value = "Primera"
$ws.Cells.Item(21, 10).Value = 250
$ws.Cells.Item(21, 11).Value = 2400
$ws.Cells.Item(21, 12).Value = 2500
$ws.Cells.Item(21, 13).Value = 2450
$ws.Cells.Item(21, 14).Value = "`$/atado 1,5 a 2 kilos"
$ws.Cells.Item(21, 16).Value = 1225
$ws.Cells.Item(21, 17).Value = 2

# Row 22 <- original row 15
$ws.Cells.Item(22, 4).Value = 44172
$ws.Cells.Item(22, 9).Value = "Primera"
$ws.Cells.Item(22, 10).Value = 200
$ws.Cells.Item(22, 11).Value = 1300
$ws.Cells.Item(22, 12).Value = 1500
$ws.Cells.Item(22, 13).Value = 1400
$ws.Cells.Item(22, 14).Value = "`$/atado 1,5 a 2 kilos"
$ws.Cells.Item(22, 16).Value = 700
$ws.Cells.Item(22, 17).Value = 2

# Row 23 <- original row 31
$ws.Cells.Item(23, 4).Value = 44936
$ws.Cells.Item(23, 9).Value = "Primera"
$ws.Cells.Item(23, 10).Value = 350
$ws.Cells.Item(23, 11).Value = 3000
$ws.Cells.Item(23, 12).Value = 3500
$ws.Cells.Item(23, 13).Value = 3357
$ws.Cells.Item(23, 14).Value = "`$/atado 1,5 a 2 kilos"
$ws.Cells.Item(23, 16).Value = 1678
$ws.Cells.Item(23, 17).Value = 2

# Row 24 <- original row 22
$ws.Cells.Item(24, 4).Value = 44616
$ws.Cells.Item(24, 9).Value = "Primera"
$ws.Cells.Item(24, 10).Value = 270
$ws.Cells.Item(24, 11).Value = 1300
$ws.Cells.Item(24, 12).Value = 1500
$ws.Cells.Item(24, 13).Value = 1400
$ws.Cells.Item(24, 14).Value = "`$/atado 1,5 a 2 kilos"
$ws.Cells.Item(24, 16).Value = 700
$ws.Cells.Item(24, 17).Value = 2

# Row 25 <- original row 7
$ws.Cells.Item(25, 4).Value = 44447
$ws.Cells.Item(25, 9).Value = "Primera"
$ws.Cells.Item(25, 10).Value = 300
$ws.Cells.Item(25, 11).Value = 900
$ws.Cells.Item(25, 12).Value = 1000
$ws.Cells.Item(25, 13).Value = 950
$ws.Cells.Item(25, 14).Value = "`$/atado 1,5 a 2 kilos"
$ws.Cells.Item(25, 16).Value = 475
$ws.Cells.Item(25, 17).Value = 2

# Row 26 <- original row 4
$ws.Cells.Item(26, 4).Value = 44438
$ws.Cells.Item(26, 9).Value = "Primera"
$ws.Cells.Item(26, 10).Value = 300
$ws.Cells.Item(26, 11).Value = 950
$ws.Cells.Item(26, 12).Value = 1000
$ws.Cells.Item(26, 13).Value = 975
$ws.Cells.Item(26, 14).Value = "`$/atado 1,5 a 2 kilos"
$ws.Cells.Item(26, 16).Value = 488
$ws.Cells.Item(26, 17).Value = 2

# Row 27 <- original row 38
$ws.Cells.Item(27, 4).Value = 44392
$ws.Cells.Item(27, 9).Value = "Primera"
$ws.Cells.Item(27, 10).Value = 250
$ws.Cells.Item(27, 11).Value = 1800
$ws.Cells.Item(27, 12).Value = 2000
$ws.Cells.Item(27, 13).Value = 1900
$ws.Cells.Item(27, 14).Value = "`$/atado 1,5 a 2 kilos"
$ws.Cells.Item(27, 16).Value = 950
$ws.Cells.Item(27, 17).Value = 2

# Row 28 <- original row 11
$ws.Cells.Item(28, 4).Value = 44917
$ws.Cells.Item(28, 9).Value = "Primera"
$ws.Cells.Item(28, 10).Value = 300
$ws.Cells.Item(28, 11).Value = 2700
$ws.Cells.Item(28, 12).Value = 3000
$ws.Cells.Item(28, 13).Value = 2850
$ws.Cells.Item(28, 14).Value = "`$/atado 1,5 a 2 kilos"
$ws.Cells.Item(28, 16).Value = 1425
$ws.Cells.Item(28, 17).Value = 2

# Row 29 <- original row 2
$ws.Cells.Item(29, 4).Value = 45008
$ws.Cells.Item(29, 9).Value = "Primera"
$ws.Cells.Item(29, 10).Value = 250
$ws.Cells.Item(29, 11).Value = 2000
$ws.Cells.Item(29, 12).Value = 2500
$ws.Cells.Item(29, 13).Value = 2200
$ws.Cells.Item(29, 14).Value = "`$/atado 1,5 a 2 kilos"
$ws.Cells.Item(29, 16).Value = 1100
$ws.Cells.Item(29, 17).Value = 2

# Row 30 <- original row 17
$ws.Cells.Item(30, 4).Value = 44365
$ws.Cells.Item(30, 9).Value = "Primera"
$ws.Cells.Item(30, 10).Value = 200
$ws.Cells.Item(30, 11).Value = 1800
$ws.Cells.Item(30, 12).Value = 2000
$ws.Cells.Item(30, 13).Value = 1900
$ws.Cells.Item(30, 14).Value = "`$/atado 1,5 a 2 kilos"
$ws.Cells.Item(30, 16).Value = 950
$ws.Cells.Item(30, 17).Value = 2

# Row 31 <- original row 28
$ws.Cells.Item(31, 4).Value = 44601
$ws.Cells.Item(31, 9).Value = "Primera"
$ws.Cells.Item(31, 10).Value = 270
$ws.Cells.Item(31, 11).Value = 2200
$ws.Cells.Item(31, 12).Value = 2500
$ws.Cells.Item(31, 13).Value = 2350
$ws.Cells.Item(31, 14).Value = "`$/atado 1,5 a 2 kilos"
$ws.Cells.Item(31, 16).Value = 1175
$ws.Cells.Item(31, 17).Value = 2

# Row 32 <- original row 24
$ws.Cells.Item(32, 4).Value = 44385
$ws.Cells.Item(32, 9).Value = "Primera"
$ws.Cells.Item(32, 10).Value = 300
$ws.Cells.Item(32, 11).Value = 2400
$ws.Cells.Item(32, 12).Value = 2500
$ws.Cells.Item(32, 13).Value = 2450
$ws.Cells.Item(32, 14).Value = "`$/atado 1,5 a 2 kilos"
$ws.Cells.Item(32, 16).Value = 1225
$ws.Cells.Item(32, 17).Value = 2

# Row 33 <- original row 42
$ws.Cells.Item(33, 4).Value = 44243
$ws.Cells.Item(33, 9).Value = "Primera"
$ws.Cells.Item(33, 10).Value = 250
$ws.Cells.Item(33, 11).Value = 1200
$ws.Cells.Item(33, 12).Value = 1300
$ws.Cells.Item(33, 13).Value = 1250
$ws.Cells.Item(33, 14).Value = "`$/atado 1,5 a 2 kilos"
$ws.Cells.Item(33, 16).Value = 625
$ws.Cells.Item(33, 17).Value = 2

# Row 34 <- original row 26
$ws.Cells.Item(34, 4).Value = 44427
$ws.Cells.Item(34, 9).Value = "Primera"
$ws.Cells.Item(34, 10).Value = 250
$ws.Cells.Item(34, 11).Value = 1300
$ws.Cells.Item(34, 12).Value = 1500
$ws.Cells.Item(34, 13).Value = 1400
$ws.Cells.Item(34, 14).Value = "`$/atado 1,5 a 2 kilos"
$ws.Cells.Item(34, 16).Value = 700
$ws.Cells.Item(34, 17).Value = 2

# Row 35 <- original row 20
$ws.Cells.Item(35, 4).Value = 44257
$ws.Cells.Item(35, 9).Value = "Primera"
$ws.Cells.Item(35, 10).Value = 500
$ws.Cells.Item(35, 11).Value = 1400
$ws.Cells.Item(35, 12).Value = 1500
$ws.Cells.Item(35, 13).Value = 1450
$ws.Cells.Item(35, 14).Value = "`$/atado 1,5 a 2 kilos"
$ws.Cells.Item(35, 16).Value = 725
$ws.Cells.Item(35, 17).Value = 2

# Row 36 <- original row 6
$ws.Cells.Item(36, 4).Value = 44544
$ws.Cells.Item(36, 9).Value = "Primera"
$ws.Cells.Item(36, 10).Value = 250
$ws.Cells.Item(36, 11).Value = 900
$ws.Cells.Item(36, 12).Value = 1000
$ws.Cells.Item(36, 13).Value = 950
$ws.Cells.Item(36, 14).Value = "`$/atado 1,5 a 2 kilos"
$ws.Cells.Item(36, 16).Value = 475
$ws.Cells.Item(36, 17).Value = 2

# Row 37 <- original row 21
$ws.Cells.Item(37, 4).Value = 44971
$ws.Cells.Item(37, 9).Value = "Primera"
$ws.Cells.Item(37, 10).Value = 350
$ws.Cells.Item(37, 11).Value = 2500
$ws.Cells.Item(37, 12).Value = 2800
$ws.Cells.Item(37, 13).Value = 2671
$ws.Cells.Item(37, 14).Value = "`$/atado 1,5 a 2 kilos"
$ws.Cells.Item(37, 16).Value = 1336
$ws.Cells.Item(37, 17).Value = 2

# Row 38 <- original row 41
$ws.Cells.Item(38, 4).Value = 45036
$ws.Cells.Item(38, 9).Value = "Segunda"
$ws.Cells.Item(38, 10).Value = 210
$ws.Cells.Item(38, 11).Value = 2300
$ws.Cells.Item(38, 12).Value = 2500
$ws.Cells.Item(38, 13).Value = 2443
$ws.Cells.Item(38, 14).Value = "`$/atado 1,5 a 2 kilos"
$ws.Cells.Item(38, 16).Value = 1222
$ws.Cells.Item(38, 17).Value = 2

# Row 39 <- original row 40
$ws.Cells.Item(39, 4).Value = 44403
$ws.Cells.Item(39, 9).Value = "Primera"
$ws.Cells.Item(39, 10).Value = 250
$ws.Cells.Item(39, 11).Value = 1800
$ws.Cells.Item(39, 12).Value = 2000
$ws.Cells.Item(39, 13).Value = 1900
$ws.Cells.Item(39, 14).Value = "`$/atado 1,5 a 2 kilos"
$ws.Cells.Item(39, 16).Value = 950
$ws.Cells.Item(39, 17).Value = 2

# Row 40 <- original row 5
$ws.Cells.Item(40, 4).Value = 44181
$ws.Cells.Item(40, 9).Value = "Primera"
$ws.Cells.Item(40, 10).Value = 200
$ws.Cells.Item(40, 11).Value = 1000
$ws.Cells.Item(40, 12).Value = 1200
$ws.Cells.Item(40, 13).Value = 1100
$ws.Cells.Item(40, 14).Value = "`$/atado"
$ws.Cells.Item(40, 16).Value = 1100
$ws.Cells.Item(40, 17).Value = 1

# Row 41 <- original row 12
$ws.Cells.Item(41, 4).Value = 44468
$ws.Cells.Item(41, 9).Value = "Primera"
$ws.Cells.Item(41, 10).Value = 300
$ws.Cells.Item(41, 11).Value = 900
$ws.Cells.Item(41, 12).Value = 1000
$ws.Cells.Item(41, 13).Value = 950
$ws.Cells.Item(41, 14).Value = "`$/atado 1,5 a 2 kilos"
$ws.Cells.Item(41, 16).Value = 475
$ws.Cells.Item(41, 17).Value = 2

# Row 42 <- original row 32
$ws.Cells.Item(42, 4).Value = 45055
$ws.Cells.Item(42, 9).Value = "Primera"
$ws.Cells.Item(42, 10).Value = 350
$ws.Cells.Item(42, 11).Value = 2300
$ws.Cells.Item(42, 12).Value = 2500
$ws.Cells.Item(42, 13).Value = 2386
$ws.Cells.Item(42, 14).Value = "`$/atado 1,5 a 2 kilos"
$ws.Cells.Item(42, 16).Value = 1193
$ws.Cells.Item(42, 17).Value = 2

# Row 44 <- original row 8
$ws.Cells.Item(44, 4).Value = 45001
$ws.Cells.Item(44, 9).Value = "Primera"
$ws.Cells.Item(44, 10).Value = 300
$ws.Cells.Item(44, 11).Value = 2000
$ws.Cells.Item(44, 12).Value = 2500
$ws.Cells.Item(44, 13).Value = 2250
$ws.Cells.Item(44, 14).Value = "`$/atado 1,5 a 2 kilos"
$ws.Cells.Item(44, 16).Value = 1125
$ws.Cells.Item(44, 17).Value = 2
